$wb = $excel.ActiveWorkbook

function SetRow($ws, $r, $h, $i, $j, $k, $l, $m, $n) {
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $j
    $ws.Range("K$r").Value = $k
    $ws.Range("L$r").Value = $l
    if ($m -ne "SKIP") {
        if ($m -eq $null) {
            $ws.Range("M$r").ClearContents()
        } else {
            $ws.Range("M$r").Value = $m
        }
    }
    if ($n -ne "SKIP") {
        if ($n -eq $null) {
            $ws.Range("N$r").ClearContents()
        } else {
            $ws.Range("N$r").Value = $n
        }
    }
}

$wsALC = $wb.Worksheets.Item("ALC")
SetRow $wsALC 5 358.41666 124.833336 592 124.833336 592 -9.833336000000003 -822
SetRow $wsALC 32 1224.1 1133 1263.1428 1133 1263.1428 -807 -1915.1428
SetRow $wsALC 40 1398.7778 1466.3334 1365 1466.3334 1365 -1291.3334 -1715
SetRow $wsALC 100 5801 7002.5 5000 7002.5 5000 -6461.5 -6082
SetRow $wsALC 131 3880 1495 4476.25 4485 13428.75 555 -23508.75
SetRow $wsALC 137 9094128 1988.5294 40007400 5965.5882 120022200 -3415.5882 -120027300

$wsARM = $wb.Worksheets.Item("ARM")
SetRow $wsARM 32 7855.234 8225.081 6326.533 8225.081 6326.533 -7938.081 -6900.533
SetRow $wsARM 45 2882 3362.4 2615.111 3362.4 2615.111 -2985.4 -3369.111

$wsBSM = $wb.Worksheets.Item("BSM")
SetRow $wsBSM 82 13967.143 13967.143 0 13967.143 0 -13584.143 $null
SetRow $wsBSM 85 13967.143 13967.143 0 13967.143 0 -12641.143 $null
SetRow $wsBSM 99 0 0 0 0 0 $null $null
SetRow $wsBSM 105 2767.3547 1526.3636 3449.9 1526.3636 3449.9 220.6364000000001 -6943.9
SetRow $wsBSM 122 58866.832 0 58866.832 0 58866.832 "SKIP" -68666.83199999999

$wsCRP = $wb.Worksheets.Item("CRP")
SetRow $wsCRP 16 1088 886.55554 1347 886.55554 1347 -599.55554 -1921
SetRow $wsCRP 58 2497.9644 1085.4117 4681 1085.4117 4681 -882.4117000000001 -5087
SetRow $wsCRP 113 1088 886.55554 1347 886.55554 1347 1283.44446 -5687
SetRow $wsCRP 136 2497.9644 1085.4117 4681 3256.2351 14043 -706.2351000000003 -19143

$wsCUL = $wb.Worksheets.Item("CUL")
SetRow $wsCUL 68 1195.8674 900.0893 1809.3334 2700.2679 5428.0002 -1889.2679 -7050.0002
SetRow $wsCUL 71 1195.8674 900.0893 1809.3334 8100.803699999999 16284.0006 -4044.803699999999 -24396.0006
SetRow $wsCUL 74 8985.571 5400 13766.333 16200 41298.999 -15139 -43420.999
SetRow $wsCUL 77 8985.571 5400 13766.333 48600 123896.997 -43296 -134504.997
SetRow $wsCUL 129 4953.1304 6488.8887 3965.8572 19466.6661 11897.5716 -14466.6661 -21897.5716
SetRow $wsCUL 131 872.9400000000001 298 903.2 894 2709.6 4146 -12789.6

$wsLTW = $wb.Worksheets.Item("LTW")
SetRow $wsLTW 46 966.6667 950 1000 950 1000 -762 -1376
SetRow $wsLTW 93 1282.2222 979.63635 1757.7142 979.63635 1757.7142 268.36365 -4253.7142
SetRow $wsLTW 141 71989.336 0 71989.336 0 71989.336 "SKIP" -82349.336

$wsWVR = $wb.Worksheets.Item("WVR")
SetRow $wsWVR 46 46276.332 0 46276.332 0 46276.332 "SKIP" -46738.332
SetRow $wsWVR 96 2534.5217 1624.5 3527.2727 1624.5 3527.2727 -251.5 -6273.2727
SetRow $wsWVR 109 0 0 0 0 0 "SKIP" $null
SetRow $wsWVR 134 46276.332 0 46276.332 0 138828.996 "SKIP" -143898.996
